$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates scraped from the crypto price refresh: Coin/Link swap for
# rows 35-36 (Monero <-> NEARProtocol reordered), plus refreshed Price (D) and
# Volume(1h) (E) text for every listed coin.
$updates = @(
    @{ Row=2; D="61.649.83"; E="  +1.12%  " },
    @{ Row=3; D="3.395.77"; E="  +0.35%  " },
    @{ Row=4; E="  -0.08%  " },
    @{ Row=5; D="577.13"; E="  +1.07%  " },
    @{ Row=6; D="142.81"; E="  +0.83%  " },
    @{ Row=7; E="  +0.02%  " },
    @{ Row=8; D="0.473"; E="  -0.35%  " },
    @{ Row=9; E="  +0.05%  " },
    @{ Row=10; D="0.123"; E="  -0.54%  " },
    @{ Row=11; D="0.386"; E="  -1.02%  " },
    @{ Row=12; D="3.974.67"; E="  +0.27%  " },
    @{ Row=13; E="  -0.25%  " },
    @{ Row=14; D="27.99"; E="  +1.20%  " },
    @{ Row=15; D="3.397.50"; E="  -0.37%  " },
    @{ Row=16; E="  +0.05%  " },
    @{ Row=17; D="61.678.21"; E="  +0.97%  " },
    @{ Row=18; D="6.13"; E="  +0.77%  " },
    @{ Row=19; D="13.65"; E="  +0.20%  " },
    @{ Row=20; D="9.15"; E="  +2.47%  " },
    @{ Row=21; D="389.51"; E="  +2.15%  " },
    @{ Row=22; D="74.76"; E="  -0.14%  " },
    @{ Row=23; D="0.548"; E="  -0.53%  " },
    @{ Row=24; E="  -0.03%  " },
    @{ Row=25; E="  -2.92%  " },
    @{ Row=26; D="0.182"; E="  +1.30%  " },
    @{ Row=27; D="0.997"; E="  -0.30%  " },
    @{ Row=28; D="7.41"; E="  +1.89%  " },
    @{ Row=29; D="7.99"; E="  +0.67%  " },
    @{ Row=30; E="  -0.55%  " },
    @{ Row=31; D="1.42"; E="  +0.79%  " },
    @{ Row=32; E="  -0.03%  " },
    @{ Row=33; D="23.37"; E="  +0.36%  " },
    @{ Row=34; D="6.94"; E="  -0.03%  " },
    @{ Row=35; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="5.11"; E="  +2.13%  " },
    @{ Row=36; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="168.23"; E="  +1.41%  " },
    @{ Row=37; D="3.428.30"; E="  +0.39%  " },
    @{ Row=38; D="1.48"; E="  +0.57%  " },
    @{ Row=39; D="0.0764"; E="  -0.40%  " },
    @{ Row=40; D="26.96"; E="  -0.29%  " },
    @{ Row=41; D="0.783"; E="  +0.51%  " },
    @{ Row=42; D="4.44"; E="  +1.45%  " },
    @{ Row=43; D="1.67"; E="  +0.46%  " },
    @{ Row=44; D="1.16"; E="  +2.47%  " },
    @{ Row=45; D="2.475.37"; E="  +0.75%  " },
    @{ Row=46; D="22.73"; E="  -1.11%  " },
    @{ Row=47; D="6.66"; E="  -0.77%  " },
    @{ Row=48; E="  -0.01%  " },
    @{ Row=49; D="0.0264"; E="  -0.83%  " },
    @{ Row=50; D="2.03"; E="  -4.50%  " },
    @{ Row=51; D="0.207"; E="  -0.75%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $ws.Range("B" + $u.Row).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C" + $u.Row).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Price column is stored as literal text (e.g. "61.649.83",
        # "577.13") even when it looks like a number, so force text
        # formatting before the write and restore the original "Normal"
        # style afterwards to avoid leaving a stray numeric style behind.
        $cell = $ws.Range("D" + $u.Row)
        $cell.Style = "Normal"
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
